$d = $word.ActiveDocument

# Locate the paragraph "Sede dello stage: {A_SEDE}" by index.
$idx = 0
$targetIdx = -1
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Sede dello stage: {A_SEDE}*") {
        $targetIdx = $idx
    }
}

if ($targetIdx -ge 0) {
    $p = $d.Paragraphs.Item($targetIdx)

    # Rename this paragraph's text to "Sede legale: {A_SEDE_LEGALE}"
    $p.Range.Text = "Sede legale: {A_SEDE_LEGALE}"

    # Insert a brand-new paragraph right after it, inheriting formatting.
    $p.Range.InsertParagraphAfter()

    # Set the text of the newly created paragraph.
    $newP = $d.Paragraphs.Item($targetIdx + 1)
    $newP.Range.Text = "Sede dello stage: {A_SEDE_SVOLGIMENTO}"
}
